# Insert a new weekly price record for "Caigua" (Agrícola del Norte S.A. de Arica)
# above the existing row 41 — this pushes the old rows 41-58 down to 42-59
# (the sheet's used range grows from A1:R58 to A1:R59), and then populate the
# newly inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("41:41").Insert()

$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 44466
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112036
$ws.Range("G41").Value = "Caigua"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 130
$ws.Range("K41").Value = 7000
$ws.Range("L41").Value = 8000
$ws.Range("M41").Value = 7500
$ws.Range("N41").Value = "$/caja 20 kilos"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 375
$ws.Range("Q41").Value = 20
$ws.Range("R41").Value = "Hortaliza"
